$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated (recalculated) figures for rows 2-5, columns B-F
$ws.Range("B2").Value = 339709981.53414017
$ws.Range("C2").Value = 404622358.8813218
$ws.Range("D2").Value = 469534736.228504
$ws.Range("E2").Value = 534447113.5756865
$ws.Range("F2").Value = 599359490.9228686

$ws.Range("B3").Value = 730929478.3785502
$ws.Range("C3").Value = 795841855.725732
$ws.Range("D3").Value = 860754233.072914
$ws.Range("E3").Value = 925666610.4200965
$ws.Range("F3").Value = 990578987.7672788

$ws.Range("B4").Value = 1513717797.8885026
$ws.Range("C4").Value = 1578630175.2356844
$ws.Range("D4").Value = 1643542552.5828667
$ws.Range("E4").Value = 1708454929.930049
$ws.Range("F4").Value = 1773367307.2772312

$ws.Range("B5").Value = 2453678911.083649
$ws.Range("C5").Value = 2518591288.430831
$ws.Range("D5").Value = 2583503665.778013
$ws.Range("E5").Value = 2648416043.1251955
$ws.Range("F5").Value = 2713328420.472378
